$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet from "Sheet1" to "_"
$ws.Name = "_"

# 2. Prepare the formatting for the five new data rows (25-29) BEFORE touching the
#    header style, so the new cells pick up the existing (soon to be freed up)
#    cellXfs entries instead of creating brand-new style records.
#    Column A keeps the date-formatted style used by the rest of column A.
$srcDate = $ws.Range("A2")
$dstDate = $ws.Range("A25:A29")
$srcDate.Copy()
$dstDate.PasteSpecial(-4122)  # xlPasteFormats

#    Columns B:G reuse the style that the header row currently has (font + full
#    box border, no number format) - this is the style that will be "freed up"
#    from the header in step 4 below.
$srcData = $ws.Range("B1:G1")
$dstData = $ws.Range("B25:G29")
$srcData.Copy()
$dstData.PasteSpecial(-4122)  # xlPasteFormats

$ws.Application.CutCopyMode = $false

# 3. Fill in the values for the five new rows.
$data = @(
    @(44218, 4.4000000000000004, 5.3, 4.3, 5.6, 4.5, 6),
    @(44219, 4.7, 4.0999999999999996, 4.5, 6.4, 4.9000000000000004, 5.8),
    @(44220, 4.7, 4.3, 4.9000000000000004, 5, 5.0999999999999996, 5.7),
    @(44221, 4.4000000000000004, 5.2, 4.8, 4.9000000000000004, 4.5999999999999996, 6.7),
    @(44222, 4.4000000000000004, 5.3, 4.8, 5.7, 4.4000000000000004, 4.5)
)

$row = 25
foreach ($values in $data) {
    $col = 1
    foreach ($value in $values) {
        $ws.Cells.Item($row, $col).Value = $value
        $col++
    }
    $row++
}

# 4. Give the header row its own border (no bottom edge, so it doesn't double up
#    against the border of the first data row anymore).
$header = $ws.Range("A1:G1")
$header.Borders.Item(9).LineStyle = -4142  # xlEdgeBottom -> xlLineStyleNone

# 5. Move the selection down to the first empty row below the new data.
$selResult = $ws.Range("A30").Select()
